$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the original row 1 (A1 = 0, bold font + border style), shifting row 2 (shared string) up to row 1
$ws.Rows(1).Delete()

# Set the new text value for A1
$newValue = @'
questions = [
    {
        "title": "Which of the following answers is an invalid template statement?@Component({\n template: `\n  &ltbutton #Answer1 (click)=\"log()\"&gtText 1&lt/button&gt\n  &ltbutton #Answer2 (click)=\"log(Math.max(1, 2))\"&gtText 2&lt/button&gt\n  &ltbutton #Answer3 (click)=\"log('1', '2')\"&gtText 3&lt/button&gt\n  &ltbutton #Answer4 (click)=\"log(Answer1)\"&gtText 4&lt/button&gt\n `\n})\nexport class DemoComponent {\n log(answer?: HTMLElement): void { }\n}",
        "ques_type": 2,
        "options": [
            "Answer1",
            "Answer2",
            "Answer3",
            "Answer4"
        ],
        "score": "Answer2"
    },
    {
        "title": "Which of the following is an incorrect syntax for adding the 'valid' class to an html element?",
        "ques_type": 2,
        "options": [
            "[ngClass]=\"{'valid': isValid}\"\t",
            "*ngClass=\"{'valid': isValid}\"",
            "[class.valid]=\"isValid\"",
            "[ngClass]=\"isValid ? 'valid' : ''\""
        ],
        "score": "*ngClass=\"{'valid': isValid}\""
    },
    {
        "title": "True or false: Given the service from the code snippet below, it is required to register the service as a provider in a specific NgModule to make it available at the root level.@Injectable({ providedIn: 'root' })\nexport class DemoService { }",
        "ques_type": 11,
        "options": [
            "true",
            "false"
        ],
        "score": "False"
    },
    {
        "title": "Consider the following definition of the DemoModule.What is the type of the forRoot() static method?@NgModule({\n imports: [CommonModule]\n})\nexport class DemoModule {\n static forRoot(): ? {\n  return {\n   ngModule: DemoModule,\n   providers: [DemoService]\n  }\n }\n}",
        "ques_type": 2,
        "options": [
            "LazyLoadedModule",
            "ModuleWithProviders",
            "ModuleForRoot",
            "Observable&ltNgModule&gt"
        ],
        "score": "ModuleWithProviders"
    }
]
'@

$ws.Range("A1").Value = $newValue

# Excel auto-adjusts row height for multi-line text; reset it back to the default since
# the target file does not carry an explicit row height / customHeight attribute.
$ws.Rows(1).AutoFit()
